$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28 (pushes existing rows 28..37 down to 29..38)
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new data record
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C28").Value = 'Arica y Parinacota'
$ws.Range("D28").Value = 44460
$ws.Range("E28").Value = 15
$ws.Range("F28").Value = 100112009
$ws.Range("G28").Value = 'Acelga'
$ws.Range("H28").Value = 'Sin especificar'
$ws.Range("I28").Value = 'Primera'
$ws.Range("J28").Value = 250
$ws.Range("K28").Value = 1400
$ws.Range("L28").Value = 1500
$ws.Range("M28").Value = 1450
$ws.Range("N28").Value = '$/atado 2,5 a 3 kilos'
$ws.Range("O28").Value = 'Región de Arica y Parinacota'
$ws.Range("P28").Value = 483
$ws.Range("Q28").Value = 3
$ws.Range("R28").Value = 'Hortaliza'
